# Applies the "Google_chrome_browser" slicing example added below the
# existing "Information_Technology" example on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 61: plain index ruler 0..20 across AH61:BB61
# ---------------------------------------------------------------------
for ($i = 0; $i -le 20; $i++) {
    $ws.Cells.Item(61, 34 + $i).Value = $i
}

# ---------------------------------------------------------------------
# Row 62: the sample string "Google_chrome_browser", one letter per
# cell AH62:BB62, highlighted with dark-red text on a yellow fill.
# ---------------------------------------------------------------------
$sample = "Google_chrome_browser"
for ($i = 0; $i -lt $sample.Length; $i++) {
    $ws.Cells.Item(62, 34 + $i).Value = $sample.Substring($i, 1)
}

# Build the highlight format once on a scratch cell, then copy it onto
# the whole row in a single paste so only one new style record is made.
$fmtSrc = $ws.Range("ZZ500")
$fmtSrc.HorizontalAlignment = -4108
$fmtSrc.Font.Color = 192
$fmtSrc.Interior.Color = 65535
$fmtSrc.Copy()
$ws.Range("AH62:BB62").PasteSpecial(-4122)
$fmtSrc.Clear()

# ---------------------------------------------------------------------
# Row 63: matching negative index ruler -21..-1 across AH63:BB63
# ---------------------------------------------------------------------
for ($i = 0; $i -le 20; $i++) {
    $ws.Cells.Item(63, 34 + $i).Value = $i - 21
}

# ---------------------------------------------------------------------
# Rows 67-75: slice-expression reference table in columns AG (expr)
# and AI (result), with AG67:AH67 merged and centred like the other
# merged caption cells on this sheet.
# ---------------------------------------------------------------------
$ws.Range("AG67:AH67").Merge()
$ws.Range("AG67:AH67").HorizontalAlignment = -4108

$ws.Cells.Item(67, 33).Value = "a[-10:10:3]"
$ws.Cells.Item(67, 35).Value = "empty"

$ws.Cells.Item(68, 33).Value = "[-5:-1:1]"
$ws.Cells.Item(68, 35).Value = "owse"

$ws.Cells.Item(69, 33).Value = "[-8:-2:10]"
$ws.Cells.Item(69, 35).Value = "_"

$ws.Cells.Item(70, 33).Value = "[:5:2]"
$ws.Cells.Item(70, 35).Value = "Gol"

$ws.Cells.Item(71, 33).Value = "[10:14:2]"
$ws.Cells.Item(71, 35).Value = "oe"

$ws.Cells.Item(72, 33).Value = "[-10:-2:3]"
$ws.Cells.Item(72, 35).Value = "mbw"

$ws.Cells.Item(73, 33).Value = "[5:12:5]"
$ws.Cells.Item(73, 35).Value = "eo"

$ws.Cells.Item(74, 33).Value = "[3:-3:3]"
$ws.Cells.Item(74, 35).Value = "g_rer"

$ws.Cells.Item(75, 33).Value = "[-5:-3:20]"
$ws.Cells.Item(75, 35).Value = "o"

# ---------------------------------------------------------------------
# Move the view roughly where the author left it.
# ---------------------------------------------------------------------
$ws.Range("AN76").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 32
